# Weekly update for "Hortaliza, Feria Lagunitas de Puerto Montt - Ajo":
# a new week's price record is inserted at the top of the data (row 47,
# right after the header block that stays frozen in rows 2-46), pushing
# every existing record down by one row. The very last existing record
# (old row 179) ends up as the new last row (180).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 47:179 down to 48:180, leaving a blank row 47 behind with the
# formatting copied from the surrounding rows (matches Excel's native
# "Insert Copied/Cut Cells" down-shift behaviour).
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44525
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112003
$ws.Cells.Item(47, 7).Value = "Ajo"
$ws.Cells.Item(47, 8).Value = "Chino"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 80
$ws.Cells.Item(47, 11).Value = 21000
$ws.Cells.Item(47, 12).Value = 21000
$ws.Cells.Item(47, 13).Value = 21000
$ws.Cells.Item(47, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(47, 15).Value = "China"
$ws.Cells.Item(47, 16).Value = 2100
$ws.Cells.Item(47, 17).Value = 10
$ws.Cells.Item(47, 18).Value = "Hortaliza"
